$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 16; this shifts existing rows 16-47 down to 17-48
# (Excel's Insert copies formatting from the row above by default, which matches
# the existing "blank row cloned" pattern the diff implies).
$ws.Rows.Item(16).Insert()

# Populate the newly inserted row 16 with the new data record.
$ws.Range("A16").Value = 5
$ws.Range("B16").Value = "Macroferia Regional de Talca"
$ws.Range("C16").Value = "Maule"
$ws.Range("D16").Value = 44497
$ws.Range("E16").Value = 7
$ws.Range("F16").Value = 100112026
$ws.Range("G16").Value = "Haba"
$ws.Range("H16").Value = "Sin especificar"
$ws.Range("I16").Value = "Primera"
$ws.Range("J16").Value = 500
$ws.Range("K16").Value = 5500
$ws.Range("L16").Value = 5500
$ws.Range("M16").Value = 5500
$ws.Range("N16").Value = "$/saco 25 kilos"
$ws.Range("O16").Value = "Región del Maule"
$ws.Range("P16").Value = 220
$ws.Range("Q16").Value = 25
$ws.Range("R16").Value = "Hortaliza"
